# Consolidate "Start Date" / "End Date" into a single "Date" column.
#
# For every data row the "Date" value (DD/MM/2024) is derived from the
# "Από DD.MM." text already present in the Availability column (C), and
# the now-redundant "End Date" column (E) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row from column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Rename the header in D1 from "Start Date" to "Date".
$ws.Range("D1").Value2 = "Date"

# Temporarily force column D to text format so the "DD/MM/2024" strings
# are stored as literal text rather than being auto-converted to date
# serial numbers, then clear the formatting again so no style residue
# is left behind on the cells.
$dataRange = $ws.Range($ws.Cells.Item(2, 4), $ws.Cells.Item($lastRow, 4))
$dataRange.NumberFormat = "@"

for ($row = 2; $row -le $lastRow; $row++) {
    $availability = $ws.Cells.Item($row, 3).Value2
    if ($availability -match '(\d{2})\.(\d{2})\.') {
        $day = $matches[1]
        $month = $matches[2]
        $ws.Cells.Item($row, 4).Value2 = "$day/$month/2024"
    }
}

$dataRange.ClearFormats()

# Remove column E ("End Date") entirely; the dimension becomes A1:D33.
$ws.Columns("E").Delete()
